# get_dash/user_expense.xlsx bug fix:
#  - insert three duplicate "index" columns (Unnamed: 0.1 / .1.1 / .1.1.1)
#    between the existing "Unnamed: 0" column and the "type" column
#  - correct a handful of mis-typed rows (transportation -> entertain, etc.)
#  - append several new rows of expense data (through row 19)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert 3 new columns at C (shifts old C,D,E -> F,G,H)
$ws.Range("C1:E1").EntireColumn.Insert()

# 2. New header cells for the inserted columns
$ws.Range("C1").Value = "Unnamed: 0.1"
$ws.Range("D1").Value = "Unnamed: 0.1.1"
$ws.Range("E1").Value = "Unnamed: 0.1.1.1"

# 3. Give the new rows (12-19) the same look as the existing "A" column
#    (bold / bordered / centered) by copying the style from row 11.
$ws.Range("A11").Copy($ws.Range("A12:A19"))

# 4. Full target data for rows 2-19
$data = @(
  @{ row=2;  A=0;  B=0;  C=0;  D=0;  E=0;     F='food';       G=10;  H='2021-08-29' },
  @{ row=3;  A=1;  B=2;  C=2;  D=2;  E=2;     F='food';       G=20;  H='2021-06-14' },
  @{ row=4;  A=2;  B=3;  C=3;  D=3;  E=3;     F='entertain';  G=80;  H='2021-08-20' },
  @{ row=5;  A=3;  B=4;  C=4;  D=4;  E=4;     F='other';      G=422; H='2021-06-14' },
  @{ row=6;  A=4;  B=5;  C=5;  D=5;  E=5;     F='clothing';   G=60.9; H='2021-08-21' },
  @{ row=7;  A=5;  B=6;  C=6;  D=6;  E=6;     F='clothing';   G=20;  H='2021-08-30' },
  @{ row=8;  A=6;  B=7;  C=7;  D=7;  E=7;     F='study';      G=20;  H='2021-08-30' },
  @{ row=9;  A=7;  B=8;  C=8;  D=8;  E=8;     F='food';       G=4.8; H='2021-08-30' },
  @{ row=10; A=8;  B=9;  C=9;  D=9;  E=9;     F='others';     G=20;  H='2021-08-30' },
  @{ row=11; A=9;  B=10; C=10; D=10; E=10;    F='food';       G=40;  H='2021-08-30' },
  @{ row=12; A=10; B=11; C=11; D=11; E=11;    F='food';       G=40;  H='2021-08-30' },
  @{ row=13; A=11; B=12; C=12; D=12; E=12;    F='others';     G=20;  H='2021-08-30' },
  @{ row=14; A=12; B=13; C=13; D=13; E=13;    F='clothing';   G=30;  H='2021-09-01' },
  @{ row=15; A=13; B=14; C=14; D=14; E=14;    F='transport';  G=300; H='2021-09-03' },
  @{ row=16; A=14; B=15; C=15; D=15; E=15;    F='clothing';   G=40;  H='2021-08-31' },
  @{ row=17; A=15; B=16; C=16; D=16; E=$null; F='entertain';  G=25;  H='2021-09-01' },
  @{ row=18; A=16; B=17; C=$null; D=$null; E=$null; F='study';    G=26;  H='2021-09-02' },
  @{ row=19; A=17; B=$null; C=$null; D=$null; E=$null; F='clothing'; G=10; H='2021-09-02' }
)

foreach ($r in $data) {
  $ws.Cells.Item($r.row, 1).Value = $r.A
  if ($null -ne $r.B) { $ws.Cells.Item($r.row, 2).Value = $r.B }
  if ($null -ne $r.C) { $ws.Cells.Item($r.row, 3).Value = $r.C }
  if ($null -ne $r.D) { $ws.Cells.Item($r.row, 4).Value = $r.D }
  if ($null -ne $r.E) { $ws.Cells.Item($r.row, 5).Value = $r.E }
  $ws.Cells.Item($r.row, 6).Value = $r.F
  $ws.Cells.Item($r.row, 7).Value = $r.G
  # The date column holds plain text like "2021-08-30"; if assigned directly
  # Excel auto-converts it into a real date serial number. Write it as a
  # formula returning the literal string, then freeze it back down to a
  # plain value so no formula (and no extra date number-format/style) is
  # left behind.
  $ws.Cells.Item($r.row, 8).Formula = "=""" + $r.H + """"
}

$dateRange = $ws.Range("H2:H19")
$dateRange.Copy()
$dateRange.PasteSpecial(-4163)
